$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns remain stored as plain text,
# matching the original inlineStr cell contents (avoids Excel auto-converting
# numeric-looking strings like "0.9998" or "30.583.21" into numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.583.21'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '1.916.40'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("D5").Value = '244.99'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = '0.4831'
$ws.Range("E7").Value = '  +1.80%  '
$ws.Range("D8").Value = '0.2897'
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '0.06794'
$ws.Range("E9").Value = '  -0.71%  '
$ws.Range("D10").Value = '112.18'
$ws.Range("E10").Value = '  +6.48%  '
$ws.Range("D11").Value = '19.54'
$ws.Range("E11").Value = '  +6.13%  '
$ws.Range("D12").Value = '1.916.54'
$ws.Range("E12").Value = '  -0.31%  '
$ws.Range("D13").Value = '0.07577'
$ws.Range("E13").Value = '  -1.43%  '
$ws.Range("D14").Value = '5.401'
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("D15").Value = '0.6720'
$ws.Range("D16").Value = '294.59'
$ws.Range("E16").Value = '  +1.52%  '
$ws.Range("D17").Value = '30.580.17'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("E18").Value = '  +0.69%  '
$ws.Range("D19").Value = '0.000007604'
$ws.Range("D20").Value = '0.9998'
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").Value = '5.517'
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("D22").Value = '2.165.69'
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("D23").Value = '0.9999'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = '6.419'
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("D25").Value = '9.484'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("D26").Value = '166.16'
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").Value = '20.28'
$ws.Range("E27").Value = '  -4.20%  '
$ws.Range("D28").Value = '2.092'
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("E29").Value = '  -0.60%  '
$ws.Range("D30").Value = '1.443'
$ws.Range("E30").Value = '  +2.92%  '
$ws.Range("D31").Value = '4.121'
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("D32").Value = '4.052'
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").Value = '0.04987'
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("D34").Value = '0.7349'
$ws.Range("E34").Value = '  +0.46%  '
$ws.Range("D35").Value = '1.143'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '2.714'
$ws.Range("E36").Value = '  -0.71%  '
$ws.Range("D37").Value = '0.02028'
$ws.Range("E37").Value = '  -2.08%  '
$ws.Range("D38").Value = '2.684'
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").Value = '2.022'
$ws.Range("E39").Value = '  -1.35%  '
$ws.Range("D40").Value = '109.44'
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("D41").Value = '0.4433'
$ws.Range("E41").Value = '  +0.67%  '
$ws.Range("D42").Value = '0.8651'
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = '5.844'
$ws.Range("E43").Value = '  -1.54%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '69.43'
$ws.Range("E45").Value = '  +2.29%  '
$ws.Range("D46").Value = '7.218'
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("D47").Value = '48.68'
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("D48").Value = '9.174'
$ws.Range("E48").Value = '  -1.71%  '
$ws.Range("D49").Value = '0.1225'
$ws.Range("E49").Value = '  -1.55%  '
$ws.Range("B50").Value = 'WOONetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D50").Value = '0.2505'
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '34.73'
$ws.Range("E51").Value = '  -0.80%  '
